$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "location" value for the user row (C2) to include ", ISB"
$ws.Range("C2").Value = "H-12 Sector, ISB"
